$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at X:Y (pushing old codelist.. columns right by 2)
$ws.Range("X1:Y1").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("X1").Value = "derived_variable"
$ws.Range("Y1").Value = "derivation_description"

# Approximate the original column widths as closely as the engine's pixel-quantized
# ColumnWidth setter allows (column X inherits ~ old column W's width, column Y gets
# its own new width)
$ws.Columns("X").ColumnWidth = 15.83
$ws.Columns("Y").ColumnWidth = 22.0

# Refresh the AutoFilter over the new full range (autofilter doesn't auto-expand
# when columns are inserted at its edge)
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ34").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Collection_AE!`$A`$1:`$AJ`$34"
    }
}

# Match the author's final selection/scroll position recorded in the workbook
$ws.Application.Goto($ws.Range("N1"))
$ws.Range("X24").Select()
